# Applies the capital-structure database refresh:
# - Removes the Sony Financial Holdings Inc. (TSE:8729) row
# - Refreshes every metric for the remaining Japan life-insurance comparables
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sony Financial Holdings Inc. (TSE:8729) dropped from the comp set;
# Excel shifts the remaining company rows up automatically.
$ws.Rows(5).Delete()

# --- Row 2: 5 ---
$ws.Range("A2").Value = "Japan"
$ws.Range("B2").Value = "'5"
$ws.Range("C2").Value = "Insurance (Life)"
$ws.Range("D2").Value = -0.008774999999999998
$ws.Range("E2").Value = 0.102
$ws.Range("F2").Value = 0.04185
$ws.Range("G2").Value = 0.06364508428564568
$ws.Range("H2").Value = 0.06364508428564568
$ws.Range("I2").Value = 0.04770714656125215
$ws.Range("J2").Value = 0.03573736891124991
$ws.Range("K2").Value = 7223.6
$ws.Range("L2").Value = 0.03125571960717035
$ws.Range("M2").Value = 2490.2482
$ws.Range("N2").Value = 0.03685876507321435
$ws.Range("O2").Value = 0.3447378315521347
$ws.Range("P2").Value = 2059.7442
$ws.Range("Q2").Value = 0.03048677139038423
$ws.Range("R2").Value = 0.2851409546486516
$ws.Range("S2").Value = 430.504
$ws.Range("T2").Value = 0.17287594063917
$ws.Range("U2").Value = 588013.7
$ws.Range("V2").Value = 8.703332795554891
$ws.Range("W2").Value = 0.03225469053525477
$ws.Range("X2").Value = 0.08453563542811299
$ws.Range("Y2").Value = -0.05228094489285822
$ws.Range("Z2").Value = -1.966555197096035
$ws.Range("AA2").Value = 0.04305778871942243
$ws.Range("AB2").Value = 0.05843486618536635
$ws.Range("AC2").Value = -0.01246171481555975
$ws.Range("AD2").Value = 296498.3
$ws.Range("AE2").Value = 836.5450375199371
$ws.Range("AF2").Value = 297334.8450375199
$ws.Range("AG2").Value = -290678.85496248
$ws.Range("AH2").Value = 0.8148465259862675
$ws.Range("AI2").Value = 0.5715841547952124
$ws.Range("AJ2").Value = 1.302809349524161
$ws.Range("AK2").Value = 4.286075193774601
$ws.Range("AL2").Value = 458.2
$ws.Range("AM2").Value = 458.2
$ws.Range("AN2").Value = 23.59876270142038
$ws.Range("AO2").Value = 24.14447839371453
$ws.Range("AP2").Value = -23.13558398338256
$ws.Range("AQ2").Value = 24.14447839371453

# --- Row 3: T&D Holdings, Inc. (TSE:8795) ---
$ws.Range("A3").Value = "Japan"
$ws.Range("B3").Value = "T&D Holdings, Inc. (TSE:8795)"
$ws.Range("C3").Value = "Insurance (Life)"
$ws.Range("D3").Value = 0.00455
$ws.Range("E3").Value = 0.102
$ws.Range("F3").Value = 0.102
$ws.Range("G3").Value = 0.1146840314025832
$ws.Range("H3").Value = 0.1146840314025832
$ws.Range("I3").Value = 0.1341084530379205
$ws.Range("J3").Value = 0.1118903046303734
$ws.Range("K3").Value = 1385.2
$ws.Range("L3").Value = 0.06484624063142225
$ws.Range("M3").Value = 246.3008
$ws.Range("N3").Value = 0.03538194564155605
$ws.Range("O3").Value = 0.1778088362691308
$ws.Range("P3").Value = 246.1968
$ws.Range("Q3").Value = 0.03536700568867436
$ws.Range("R3").Value = 0.1777337568582154
$ws.Range("S3").Value = 0.1040000000000134
$ws.Range("T3").Value = 0.0004222479179930127
$ws.Range("U3").Value = 9117.0
$ws.Range("V3").Value = 1.309687984830202
$ws.Range("W3").Value = 0.1221419816769392
$ws.Range("X3").Value = 0.06938471719349035
$ws.Range("Y3").Value = 0.05275726448344888
$ws.Range("Z3").Value = 3.977025133105832
$ws.Range("AA3").Value = 0.444990553665863
$ws.Range("AB3").Value = 0.06003035270110289
$ws.Range("AC3").Value = 0.3849602009647601
$ws.Range("AD3").Value = 1535.9
$ws.Range("AE3").Value = 0.0755106053421718
$ws.Range("AF3").Value = 1535.975510605342
$ws.Range("AG3").Value = -7581.024489394658
$ws.Range("AH3").Value = 0.1807630675261783
$ws.Range("AI3").Value = 0.1043463360043441
$ws.Range("AJ3").Value = 12.23092120286913
$ws.Range("AK3").Value = -1.35303544965443
$ws.Range("AL3").Value = 12.8
$ws.Range("AM3").Value = 12.8
$ws.Range("AN3").Value = 0.5038635333237541
$ws.Range("AO3").Value = 223.8046875
$ws.Range("AP3").Value = -2.487012035575429
$ws.Range("AQ3").Value = 223.8046875

# --- Row 4: Dai-ichi Life Holdings, Inc. (TSE:8750) ---
$ws.Range("A4").Value = "Japan"
$ws.Range("B4").Value = "Dai-ichi Life Holdings, Inc. (TSE:8750)"
$ws.Range("C4").Value = "Insurance (Life)"
$ws.Range("D4").Value = -0.0221
$ws.Range("E4").Value = -0.311
$ws.Range("F4").Value = 1.046
$ws.Range("G4").Value = 0.142368299566432
$ws.Range("H4").Value = 0.142368299566432
$ws.Range("I4").Value = 0.08491085846525136
$ws.Range("J4").Value = 0.04245542923262568
$ws.Range("K4").Value = 227.2
$ws.Range("L4").Value = 0.003604341073465418
$ws.Range("M4").Value = 1083.1976
$ws.Range("N4").Value = 0.06488622123183457
$ws.Range("O4").Value = 4.767595070422535
$ws.Range("P4").Value = 652.7976
$ws.Range("Q4").Value = 0.03910419437156309
$ws.Range("R4").Value = 2.873228873239437
$ws.Range("S4").Value = 430.4
$ws.Range("T4").Value = 0.3973420915999075
$ws.Range("U4").Value = 13332.6
$ws.Range("V4").Value = 0.7986557883765231
$ws.Range("W4").Value = 0.005677473524281694
$ws.Range("X4").Value = 0.08453563542811299
$ws.Range("Y4").Value = -0.0788581619038313
$ws.Range("Z4").Value = 1.603698553606659
$ws.Range("AA4").Value = 0.06808571045311167
$ws.Range("AB4").Value = 0.05843486618536635
$ws.Range("AC4").Value = 0.009650844267745318
$ws.Range("AD4").Value = 10610.9
$ws.Range("AE4").Value = 68.6777277851738
$ws.Range("AF4").Value = 10679.57772778517
$ws.Range("AG4").Value = -2653.022272214826
$ws.Range("AH4").Value = 0.3901446812296364
$ws.Range("AI4").Value = 0.2057370183148601
$ws.Range("AJ4").Value = -0.1889512335890614
$ws.Range("AK4").Value = -0.06877341279363366
$ws.Range("AL4").Value = 424.6
$ws.Range("AM4").Value = 424.6
$ws.Range("AN4").Value = 1.773537916395059
$ws.Range("AO4").Value = 12.6073951954781
$ws.Range("AP4").Value = -0.4434341660757871
$ws.Range("AQ4").Value = 12.6073951954781

# --- Row 5: Japan Post Insurance Co., Ltd. (TSE:7181) ---
$ws.Range("A5").Value = "Japan"
$ws.Range("B5").Value = "Japan Post Insurance Co., Ltd. (TSE:7181)"
$ws.Range("C5").Value = "Insurance (Life)"
$ws.Range("D5").Value = -0.109
$ws.Range("E5").Value = 0.163
$ws.Range("F5").Value = -0.0183
$ws.Range("G5").Value = 0.08840902825545263
$ws.Range("H5").Value = 0.08840902825545263
$ws.Range("I5").Value = 0.07619155049001346
$ws.Range("J5").Value = 0.05190081119954614
$ws.Range("K5").Value = 1593.3
$ws.Range("L5").Value = 0.04264413439071588
$ws.Range("M5").Value = 202.464
$ws.Range("N5").Value = 0.01759149202377228
$ws.Range("O5").Value = 0.1270721144793824
$ws.Range("P5").Value = 202.464
$ws.Range("Q5").Value = 0.01759149202377228
$ws.Range("R5").Value = 0.1270721144793824
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 12731.7
$ws.Range("V5").Value = 1.106219372328224
$ws.Range("W5").Value = 0.0768822470673956
$ws.Range("X5").Value = 0.1843337013105319
$ws.Range("Y5").Value = -0.1074514542431363
$ws.Range("Z5").Value = 0.8296168734988666
$ws.Range("AA5").Value = 0.04305778871942243
$ws.Range("AB5").Value = 0.05551950353498218
$ws.Range("AC5").Value = -0.01246171481555975
$ws.Range("AD5").Value = 38770.4
$ws.Range("AE5").Value = 363.3897825338718
$ws.Range("AF5").Value = 39133.78978253387
$ws.Range("AG5").Value = 26402.08978253387
$ws.Range("AH5").Value = 0.7727385359864876
$ws.Range("AI5").Value = 0.623933511331169
$ws.Range("AJ5").Value = 0.6964176089492369
$ws.Range("AK5").Value = 0.5281538721994697
$ws.Range("AL5").Value = 20.8
$ws.Range("AM5").Value = 20.8
$ws.Range("AN5").Value = 11.11759814182892
$ws.Range("AO5").Value = 137.7884615384615
$ws.Range("AP5").Value = 7.570925868876745
$ws.Range("AQ5").Value = 137.7884615384615

# --- Row 6: Japan Post Holdings Co., Ltd. (TSE:6178) ---
$ws.Range("A6").Value = "Japan"
$ws.Range("B6").Value = "Japan Post Holdings Co., Ltd. (TSE:6178)"
$ws.Range("C6").Value = "Insurance (Life)"
$ws.Range("D6").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = -0.0383
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = -0.0001582531324241
$ws.Range("J6").Value = -0.0001155214202878325
$ws.Range("K6").Value = 4039.8
$ws.Range("L6").Value = 0.03699630110252604
$ws.Range("M6").Value = 958.2858
$ws.Range("N6").Value = 0.03046955542978512
$ws.Range("O6").Value = 0.2372111985741868
$ws.Range("P6").Value = 958.2858
$ws.Range("Q6").Value = 0.03046955542978512
$ws.Range("R6").Value = 0.2372111985741868
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 552812.1
$ws.Range("V6").Value = 17.57715592071375
$ws.Range("W6").Value = 0.03225469053525477
$ws.Range("X6").Value = 0.3441673451806269
$ws.Range("Y6").Value = -0.3119126546453722
$ws.Range("Z6").Value = -0.5266844763894972
$ws.Range("AA6").Value = 0.00006084333875606808
$ws.Range("AB6").Value = 0.054651587310666
$ws.Range("AC6").Value = -0.05459074397190993
$ws.Range("AD6").Value = 245581.1
$ws.Range("AE6").Value = 404.4020165955494
$ws.Range("AF6").Value = 245985.5020165956
$ws.Range("AG6").Value = -306826.5979834044
$ws.Range("AH6").Value = 0.8866384015223848
$ws.Range("AI6").Value = 0.629637130371933
$ws.Range("AJ6").Value = 1.114209663261557
$ws.Range("AK6").Value = 1.892427196284459
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
$ws.Range("AN6").Value = 3861.338050314465
$ws.Range("AO6").ClearContents()
$ws.Range("AP6").Value = -4824.317578355415
$ws.Range("AQ6").ClearContents()

# --- Row 7: Lifenet Insurance Company (TSE:7157) ---
$ws.Range("A7").Value = "Japan"
$ws.Range("B7").Value = "Lifenet Insurance Company (TSE:7157)"
$ws.Range("C7").Value = "Insurance (Life)"
$ws.Range("D7").Value = 0.137
$ws.Range("E7").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("G7").Value = -0.1131363922061596
$ws.Range("H7").Value = -0.1131363922061596
$ws.Range("I7").Value = -0.1307353865493401
$ws.Range("J7").Value = -0.1307353865493401
$ws.Range("K7").Value = -21.9
$ws.Range("L7").Value = -0.1376492771841609
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").ClearContents()
$ws.Range("U7").Value = 20.3
$ws.Range("V7").Value = 0.02143385070214339
$ws.Range("W7").Value = -0.2198795180722891
$ws.Range("X7").Value = 0.0614077660836822
$ws.Range("Y7").Value = -0.2812872841559714
$ws.Range("Z7").Value = 1.774481374079857
$ws.Range("AA7").Value = -0.2319875083649342
$ws.Range("AB7").Value = 0.0614077660836822
$ws.Range("AC7").Value = -0.2933952744486164
$ws.Range("AD7").Value = 0
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 0
$ws.Range("AG7").Value = -20.3
$ws.Range("AH7").Value = 0
$ws.Range("AI7").Value = 0
$ws.Range("AJ7").Value = -0.02190332326283988
$ws.Range("AK7").Value = -0.1391363947909527
$ws.Range("AL7").Value = 0
$ws.Range("AM7").Value = 0
$ws.Range("AN7").Value = 0
$ws.Range("AO7").ClearContents()
$ws.Range("AP7").Value = 1.134078212290503
$ws.Range("AQ7").ClearContents()

